$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.324806928634644
$ws.Range("B3").Value = 1.238573789596558
$ws.Range("B4").Value = 4.952423810958862
$ws.Range("B5").Value = 5.056061983108521
$ws.Range("B6").Value = 5.972718954086304
